$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion-of-the-day text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 14.95 = 63159.94 pesos`n✅ 63159.94 pesos = 14.87 = 972.61 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the rate cells in N10/O10/N12/O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 66.90000000000001
$ws2.Range("O10").Value = 4225.4
$ws2.Range("N12").Value = 4247
$ws2.Range("O12").Value = 65.40000000000001
